# Update "想去人数" (F column) figures for both the "展览" sheet and the
# aggregated "全部类型" sheet, reflecting refreshed counts from the source site.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (row -> new F value)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1439
$ws1.Range("F7").Value = 2145
$ws1.Range("F11").Value = 4783
$ws1.Range("F14").Value = 298
$ws1.Range("F15").Value = 221
$ws1.Range("F20").Value = 111
$ws1.Range("F21").Value = 3710
$ws1.Range("F22").Value = 595
$ws1.Range("F23").Value = 611
$ws1.Range("F30").Value = 79
$ws1.Range("F34").Value = 840
$ws1.Range("F35").Value = 2317

# Sheet "全部类型" (row -> new F value)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 1439
$ws4.Range("F7").Value = 2145
$ws4.Range("F11").Value = 4783
$ws4.Range("F14").Value = 298
$ws4.Range("F15").Value = 221
$ws4.Range("F20").Value = 111
$ws4.Range("F21").Value = 3710
$ws4.Range("F22").Value = 595
$ws4.Range("F23").Value = 611
$ws4.Range("F30").Value = 79
$ws4.Range("F35").Value = 840
$ws4.Range("F36").Value = 2317
